# Applies the 30/12/2025 21:49 re-sync update:
#  - refreshed "Atualizado"/sync timestamps
#  - Transporte, Assinaturas, Compras, Lazer actuals changed (PicPay import:
#    health insurance -> Saude bucket shuffled totals around; utilities -> Casa)
#  - Gastos Variaveis / Obra summary values recomputed
#  - Mensal (hidden col M), Categorias and Dados mirrors updated to match

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Dashboard
# ---------------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Cells.Item(2, 1).Value = "Atualizado: 30/12/2025 21:49"

# Resumo do Mes
$dash.Cells.Item(7, 3).Value = 31101.39000000001
$dash.Cells.Item(7, 4).Value = "'48%"

$dash.Cells.Item(9, 3).Value = 14398.60999999999

# Gastos por categoria
$dash.Cells.Item(15, 3).Value = 7882.21
$dash.Cells.Item(15, 4).Value = 225

$dash.Cells.Item(17, 3).Value = 3492.91
$dash.Cells.Item(17, 4).Value = 83

$dash.Cells.Item(18, 3).Value = 1967.4
$dash.Cells.Item(18, 4).Value = 51

$dash.Cells.Item(19, 3).Value = 2551.24
$dash.Cells.Item(19, 4).Value = 196

# ---------------------------------------------------------------------------
# Mensal (hidden helper column M holds per-category running totals)
# ---------------------------------------------------------------------------
$mensal = $wb.Worksheets.Item("Mensal")

$mensal.Cells.Item(5, 13).Value = 7882.21
$mensal.Cells.Item(7, 13).Value = 3492.91
$mensal.Cells.Item(8, 13).Value = 1967.4
$mensal.Cells.Item(9, 13).Value = 2551.24

# ---------------------------------------------------------------------------
# Categorias
# ---------------------------------------------------------------------------
$cats = $wb.Worksheets.Item("Categorias")

$cats.Cells.Item(5, 3).Value = 7882.21
$cats.Cells.Item(5, 4).Value = -4382.21
$cats.Cells.Item(5, 5).Value = 2.25206

$cats.Cells.Item(7, 3).Value = 3492.91
$cats.Cells.Item(7, 4).Value = 707.0900000000001
$cats.Cells.Item(7, 5).Value = 0.831645238095238

$cats.Cells.Item(8, 3).Value = 1967.4
$cats.Cells.Item(8, 4).Value = 1832.6
$cats.Cells.Item(8, 5).Value = 0.5177368421052632

$cats.Cells.Item(9, 3).Value = 2551.24
$cats.Cells.Item(9, 4).Value = -1251.24
$cats.Cells.Item(9, 5).Value = 1.962492307692308

# ---------------------------------------------------------------------------
# Dados
# ---------------------------------------------------------------------------
$dados = $wb.Worksheets.Item("Dados")

$dados.Cells.Item(3, 2).Value = "2025-12-30T21:49:28.626440"

$dados.Cells.Item(9, 4).Value = 7882.21
$dados.Cells.Item(11, 4).Value = 3492.91
$dados.Cells.Item(12, 4).Value = 1967.4
$dados.Cells.Item(13, 4).Value = 2551.24
